$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-21 04:26:58"

foreach ($sheetName in @("Главные", "Линейные")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA$row").Value = $newTimestamp
    }
}
